# The edit splits the empty "_GoBack" bookmark paragraph into two new
# paragraphs of instructions ("cp .env.example .env" and "php artisan
# key:generate"), preserving the _GoBack bookmark around "php artisan
# key:generate". Built as a literal WordprocessingML fragment (runs +
# w:proofErr spell/grammar markers exactly as authored) and dropped in via
# Range.InsertXML so the markup matches byte-for-byte instead of relying on
# Word's live auto-proofing to regenerate the proofErr wrappers.

$d = $word.ActiveDocument

# Anchor on the (hidden) "_GoBack" bookmark rather than a hard-coded
# paragraph index, then expand to the whole paragraph (incl. paragraph mark)
# that currently holds only that bookmark.
$bm = $d.Bookmarks.Item("_GoBack")
$target = $bm.Range.Duplicate
$target.Expand(4) | Out-Null

$newParagraphsXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00241DD3" w:rsidRDefault="00241DD3" w:rsidP="00076F5F"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Run “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>cp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>env.example</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>env</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>” to create .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>env</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> file</w:t></w:r></w:p><w:p w:rsidR="00241DD3" w:rsidRDefault="00241DD3" w:rsidP="00076F5F"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Run “</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>php</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> artisan </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>key</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:generate</w:t></w:r><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>” to generate system key</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($newParagraphsXml)
